$d = $word.ActiveDocument
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" mc:Ignorable="w14 w15 wp14"><w:body><w:p w:rsidR="006042BC" w:rsidRPr="00704A78" w:rsidRDefault="00704A78" w:rsidP="00704A78"><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:u w:val="single"/><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r w:rsidRPr="00704A78"><w:rPr><w:b/><w:u w:val="single"/><w:lang w:val="en-GB"/></w:rPr><w:t>Experiments for the hydraulic bench</w:t></w:r></w:p><w:p w:rsidR="00704A78" w:rsidRDefault="00704A78"><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr></w:p><w:p w:rsidR="00704A78" w:rsidRPr="00847859" w:rsidRDefault="00704A78"><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r w:rsidRPr="00847859"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>Experiment n°1: head losses</w:t></w:r></w:p><w:p w:rsidR="00704A78" w:rsidRDefault="00704A78"><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r w:rsidRPr="00847859"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">The purpose of this experiment is to show the effect of head losses on the mass flow ratio. </w:t></w:r><w:r w:rsidR="00084D3B" w:rsidRPr="00847859"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>For the experiment</w:t></w:r><w:r w:rsidRPr="00847859"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>, student have to open a single pipe and measure the head losses for different mass flow rate</w:t></w:r></w:p><w:p w:rsidR="00084D3B" w:rsidRDefault="00084D3B"><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>Experiment n°2: laminar and turbulent flow for the mass flow</w:t></w:r></w:p><w:p w:rsidR="00084D3B" w:rsidRDefault="00084D3B"><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>The purpose of this experiment is to show the effect of the flow regime on the mass flow.</w:t></w:r><w:r w:rsidR="00822028"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> Students have to use a least two pipe with different diameters in order to have in one pipe a laminar flow and in the other a turbulent flow (a third pipe for a transition flow can be added), students then should measure </w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">the pressure loss in each case </w:t></w:r></w:p><w:p w:rsidR="00822028" w:rsidRPr="00F023FC" w:rsidRDefault="00822028"><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r w:rsidRPr="00F023FC"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">Experiment n°3: flowmeter comparison </w:t></w:r></w:p><w:p w:rsidR="00822028" w:rsidRDefault="00822028"><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r w:rsidRPr="00F023FC"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">The purpose of this one is to compare different flowmeter, head losses will be measure for each flowmeter (rotameter, venturi tube, orifice plate) </w:t></w:r></w:p><w:p w:rsidR="00613E20" w:rsidRDefault="00822028" w:rsidP="00613E20"><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">Experiment n°4: </w:t></w:r><w:r w:rsidR="00613E20"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>head losses comparison experiment/simulation</w:t></w:r></w:p><w:p w:rsidR="00613E20" w:rsidRDefault="00893579" w:rsidP="00613E20"><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>Compare each la</w:t></w:r><w:r w:rsidR="00613E20"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>w</w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>s</w:t></w:r><w:r w:rsidR="00613E20"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> for head losses to the experiment.</w:t></w:r></w:p><w:p w:rsidR="00613E20" w:rsidRPr="00BC6ED2" w:rsidRDefault="00613E20"><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r w:rsidRPr="00BC6ED2"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>Experiment n°5:</w:t></w:r><w:r w:rsidR="00A6566B" w:rsidRPr="00BC6ED2"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> head losses comparison of different technical solution for pipe</w:t></w:r></w:p><w:p w:rsidR="00A6566B" w:rsidRDefault="00A6566B"><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r w:rsidRPr="00BC6ED2"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>On the circuit pipe, compare each section (the u like, the v like and the n like sections) conclude on the efficiency of each part.</w:t></w:r></w:p><w:p w:rsidR="00A6566B" w:rsidRDefault="00A6566B" w:rsidP="00B4020C"><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r w:rsidRPr="00787828"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">Experiment n°6: </w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>balancing of a</w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> hydraulic circuit.</w:t></w:r></w:p><w:p w:rsidR="00A6566B" w:rsidRDefault="00A6566B"><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>The objective of this experiment is to simulate a hydraulic circuit in a building, students will have to change the valve aperture in two pipes in order to have the same pressure loss in each pipe.</w:t></w:r></w:p><w:p w:rsidR="00A6566B" w:rsidRDefault="00A6566B"><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr></w:p><w:p w:rsidR="00B4020C" w:rsidRDefault="00AF4C47"><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">Experiment n°7: head losses for laminar and turbulent flow </w:t></w:r></w:p><w:p w:rsidR="00AF4C47" w:rsidRDefault="00B4020C"><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>Comparison between head losses in a laminar flow and in a turbulent flow.</w:t></w:r></w:p><w:p w:rsidR="00AF4C47" w:rsidRDefault="00AF4C47"><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>Experiment n°8</w:t></w:r><w:r w:rsidR="008B015E"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>:</w:t></w:r><w:r w:rsidR="00985608"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> check of the venturi effect</w:t></w:r></w:p><w:p w:rsidR="00985608" w:rsidRDefault="00985608"><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>Check the theory of the venturi effect and Bernoulli equation with the venturi tube equipped on the bench and the orifice plate</w:t></w:r><w:r w:rsidR="00443DD6"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> also equipped on</w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>.</w:t></w:r></w:p><w:p w:rsidR="00985608" w:rsidRDefault="00985608"><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r w:rsidRPr="001D7292"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">Experiment n°9: </w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">experimental determination of the pump power </w:t></w:r></w:p><w:p w:rsidR="00985608" w:rsidRDefault="00985608"><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>Experiment n°10:</w:t></w:r><w:r w:rsidR="00673A07"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> evaluation of the critical Reynold’s number</w:t></w:r></w:p><w:p w:rsidR="00904F07" w:rsidRDefault="00904F07"><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>The aim of this experiment is to find the critical Reynold’s number for different diameter of pipe.</w:t></w:r></w:p><w:p w:rsidR="00904F07" w:rsidRDefault="00904F07"><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr></w:p><w:p w:rsidR="00904F07" w:rsidRDefault="00904F07"><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr></w:p><w:p w:rsidR="00904F07" w:rsidRDefault="00904F07"><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:sectPr w:rsidR="00A6566B" w:rsidRPr="00704A78"><w:pgSz w:w="11906" w:h="16838"/><w:pgMar w:top="1417" w:right="1417" w:bottom="1417" w:left="1417" w:header="708" w:footer="708" w:gutter="0"/><w:cols w:space="708"/><w:docGrid w:linePitch="360"/></w:sectPr></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Content.InsertXML($xml)
Write-Output "Applied full-body XML replacement"
